$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles) from row 18 down to the new rows 19-22
# so the new cells pick up the same style indices as the existing rows.
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B22").PasteSpecial(-4122)

$values = @(
  @("CMAU9099877", ";DOLVA00063742"),
  @("CAIU8312898", ";DOLVA00063742"),
  @("TEMU6437272", ";DOLVA00063742"),
  @("FCIU9169820", ";DOLVA00063742"),
  @("BMOU6720890", ";DOLVA00063742"),
  @("TRLU4885536", ";DOLVA00063742"),
  @("CMAU7128390", ";DOLVA00063742"),
  @("TGHU6321448", ";DOLVA00063742"),
  @("TGHU9515160", ";DOLVA00063742"),
  @("TCNU3409791", ";DOLVA00063742"),
  @("TCLU6716222", ";DOLVA00063742"),
  @("GESU4748994", ";DOLVA00063818"),
  @("CMAU7813497", ";DOLVA00063855"),
  @("TLLU4566625", ";DOLVA00063817"),
  @("APZU4718696", ";DOLVA00063878"),
  @("CAIU7117840", ";DOLVA00063820"),
  @("CMAU9010489", ";DOLVA00063879"),
  @("ECMU9757564", ";DOLVA00063782"),
  @("TGBU5575981", ";DOLVA00063852"),
  @("MRKU4918940", ";704851"),
  @("MRKU6111026", ";869034"),
  @("APHU6741856", ";DOLVA00063900")
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i][0]
  $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Extend the duplicate-values conditional formatting to cover the new rows
$fc = $ws.Range("A1:A18").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:A22"))

$ws.Range("A1:B22").Select()
